# Updates the cryptos price/volume table (and a couple of row re-orderings /
# a new coin insertion) to match the latest scrape, per the commit:
# "Updated cryptos list on Tue Apr  9 17:59:43 UTC 2024 with GitHub Actions".
#
# All of B/C/D/E are plain text cells (t="inlineStr") in the source workbook,
# and several "numbers" (e.g. "1.00", "0.621", "68.827.12") must stay literal
# text rather than being auto-coerced to doubles by Excel. We force text entry
# by switching the cell to the "@" (Text) number format before assigning the
# value, then clear the format again afterwards so the cell keeps its original
# (default) style - only the stored value/type changes, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextCell "D2" "68.827.12"
Set-TextCell "E2" "  -3.99%  "

# Row 3
Set-TextCell "D3" "3.500.48"
Set-TextCell "E3" "  -4.47%  "

# Row 4
Set-TextCell "E4" "  +0.02%  "

# Row 5
Set-TextCell "D5" "578.42"
Set-TextCell "E5" "  -1.61%  "

# Row 6
Set-TextCell "D6" "175.57"
Set-TextCell "E6" "  -2.52%  "

# Row 7
Set-TextCell "D7" "0.621"
Set-TextCell "E7" "  -0.30%  "

# Row 8
Set-TextCell "D8" "3.491.94"
Set-TextCell "E8" "  -4.49%  "

# Row 9
Set-TextCell "E9" "  +0.05%  "

# Row 10
Set-TextCell "E10" "  -7.02%  "

# Row 11
Set-TextCell "D11" "6.56"
Set-TextCell "E11" "  +5.57%  "

# Row 12
Set-TextCell "D12" "0.602"
Set-TextCell "E12" "  -1.68%  "

# Row 13
Set-TextCell "D13" "47.19"
Set-TextCell "E13" "  -5.33%  "

# Row 14
Set-TextCell "E14" "  -3.50%  "

# Row 15
Set-TextCell "D15" "674.79"
Set-TextCell "E15" "  -1.31%  "

# Row 16
Set-TextCell "D16" "8.88"
Set-TextCell "E16" "  -1.41%  "

# Row 17
Set-TextCell "D17" "4.057.91"
Set-TextCell "E17" "  -4.40%  "

# Row 18: 'WrappedEther' -> 'WrappedBTC'
Set-TextCell "B18" "WrappedBTC"
Set-TextCell "C18" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D18" "68.824.21"
Set-TextCell "E18" "  -4.08%  "

# Row 19: 'WrappedBTC' -> 'WrappedEther'
Set-TextCell "B19" "WrappedEther"
Set-TextCell "C19" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D19" "3.499.94"
Set-TextCell "E19" "  -4.25%  "

# Row 20
Set-TextCell "E20" "  -1.73%  "

# Row 21
Set-TextCell "D21" "17.55"

# Row 22
Set-TextCell "D22" "11.16"
Set-TextCell "E22" "  -4.20%  "

# Row 23
Set-TextCell "D23" "0.902"
Set-TextCell "E23" "  -4.28%  "

# Row 24
Set-TextCell "D24" "16.34"
Set-TextCell "E24" "  -8.51%  "

# Row 25
Set-TextCell "D25" "98.12"
Set-TextCell "E25" "  -5.20%  "

# Row 26
Set-TextCell "D26" "3.85"
Set-TextCell "E26" "  -4.40%  "

# Row 27: 'Dai' -> 'LEO'
Set-TextCell "B27" "LEO"
Set-TextCell "C27" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D27" "5.80"
Set-TextCell "E27" "  -1.06%  "

# Row 28: 'ImmutableX' -> 'Dai'
Set-TextCell "B28" "Dai"
Set-TextCell "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D28" "1.00"
Set-TextCell "E28" "  +0.08%  "

# Row 29: 'RenderToken' -> 'ImmutableX'
Set-TextCell "B29" "ImmutableX"
Set-TextCell "C29" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D29" "2.65"
Set-TextCell "E29" "  -6.80%  "

# Row 30: 'EthereumClassic' -> 'RenderToken'
Set-TextCell "B30" "RenderToken"
Set-TextCell "C30" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D30" "9.41"
Set-TextCell "E30" "  -7.74%  "

# Row 31: 'Filecoin' -> 'EthereumClassic'
Set-TextCell "B31" "EthereumClassic"
Set-TextCell "C31" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D31" "32.92"
Set-TextCell "E31" "  -6.84%  "

# Row 32: 'Stacks' -> 'Filecoin'
Set-TextCell "B32" "Filecoin"
Set-TextCell "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D32" "8.75"
Set-TextCell "E32" "  -4.98%  "

# Row 33: 'NEARProtocol' -> 'Stacks'
Set-TextCell "B33" "Stacks"
Set-TextCell "C33" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D33" "3.20"
Set-TextCell "E33" "  -7.62%  "

# Row 34: 'Mantle' -> 'NEARProtocol'
Set-TextCell "B34" "NEARProtocol"
Set-TextCell "C34" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D34" "7.40"
Set-TextCell "E34" "  -0.60%  "

# Row 35: 'Bittensor' -> 'Mantle'
Set-TextCell "B35" "Mantle"
Set-TextCell "C35" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D35" "1.36"
Set-TextCell "E35" "  -5.97%  "

# Row 36: 'dogwifhat' -> 'Bittensor'
Set-TextCell "B36" "Bittensor"
Set-TextCell "C36" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D36" "567.47"
Set-TextCell "E36" "  -1.82%  "

# Row 37: 'Cosmos' -> 'dogwifhat'
Set-TextCell "B37" "dogwifhat"
Set-TextCell "C37" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D37" "3.63"
Set-TextCell "E37" "  -14.52%  "

# Row 38: 'Hedera' -> 'Cosmos'
Set-TextCell "B38" "Cosmos"
Set-TextCell "C38" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D38" "10.93"
Set-TextCell "E38" "  -3.65%  "

# Row 39: 'OKB' -> 'Hedera'
Set-TextCell "B39" "Hedera"
Set-TextCell "C39" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D39" "0.106"
Set-TextCell "E39" "  -3.18%  "

# Row 40: 'FirstDigitalUSD' -> 'OKB'
Set-TextCell "B40" "OKB"
Set-TextCell "C40" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D40" "56.76"
Set-TextCell "E40" "  -5.70%  "

# Row 41: 'VeChain' -> 'FirstDigitalUSD'
Set-TextCell "B41" "FirstDigitalUSD"
Set-TextCell "C41" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D41" "1.00"
Set-TextCell "E41" "  +0.12%  "

# Row 42: 'Kaspa' -> 'VeChain'
Set-TextCell "B42" "VeChain"
Set-TextCell "C42" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D42" "0.0440"
Set-TextCell "E42" "  -4.81%  "

# Row 43: 'TheGraph' -> 'Kaspa'
Set-TextCell "B43" "Kaspa"
Set-TextCell "C43" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D43" "0.137"
Set-TextCell "E43" "  -4.24%  "

# Row 44: 'Maker' -> 'TheGraph'
Set-TextCell "B44" "TheGraph"
Set-TextCell "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D44" "0.336"
Set-TextCell "E44" "  -3.06%  "

# Row 45: 'InjectiveProtocol' -> 'Maker'
Set-TextCell "B45" "Maker"
Set-TextCell "C45" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D45" "3.416.73"
Set-TextCell "E45" "  -8.68%  "

# Row 46: 'PEPE' -> 'InjectiveProtocol'
Set-TextCell "B46" "InjectiveProtocol"
Set-TextCell "C46" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D46" "33.41"
Set-TextCell "E46" "  -6.17%  "

# Row 47: 'ThetaToken' -> 'PEPE'
Set-TextCell "B47" "PEPE"
Set-TextCell "C47" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D47" "0.0₃0702"
Set-TextCell "E47" "  -8.37%  "

# Row 48: 'Fetch.AI' -> 'ThetaToken'
Set-TextCell "B48" "ThetaToken"
Set-TextCell "C48" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell "D48" "2.87"
Set-TextCell "E48" "  +2.37%  "

# Row 49: 'Stellar' -> 'Fetch.AI'
Set-TextCell "B49" "Fetch.AI"
Set-TextCell "C49" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D49" "2.60"
Set-TextCell "E49" "  -7.16%  "

# Row 50: 'Monero' -> 'Stellar'
Set-TextCell "B50" "Stellar"
Set-TextCell "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D50" "0.133"
Set-TextCell "E50" "  -0.68%  "

# Row 51: 'Cronos' -> 'Monero'
Set-TextCell "B51" "Monero"
Set-TextCell "C51" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D51" "133.97"
Set-TextCell "E51" "  +0.06%  "

